$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws 2 4 '26.010.90'
Set-TextCell $ws 2 5 '  +0.21%  '

# Row 3
Set-TextCell $ws 3 4 '1.641.56'
Set-TextCell $ws 3 5 '  -0.29%  '

# Row 4
Set-TextCell $ws 4 5 '  -0.68%  '

# Row 5
Set-TextCell $ws 5 4 '215.14'
Set-TextCell $ws 5 5 '  -0.38%  '

# Row 6
Set-TextCell $ws 6 4 '0.5090'
Set-TextCell $ws 6 5 '  -0.31%  '

# Row 7
Set-TextCell $ws 7 5 '  -0.43%  '

# Row 8
Set-TextCell $ws 8 4 '0.2584'
Set-TextCell $ws 8 5 '  +0.11%  '

# Row 9
Set-TextCell $ws 9 4 '0.06365'
Set-TextCell $ws 9 5 '  -0.91%  '

# Row 10
Set-TextCell $ws 10 4 '19.90'
Set-TextCell $ws 10 5 '  +1.18%  '

# Row 11
Set-TextCell $ws 11 4 '0.07758'
Set-TextCell $ws 11 5 '  -0.36%  '

# Row 12
Set-TextCell $ws 12 4 '4.283'
Set-TextCell $ws 12 5 '  -0.86%  '

# Row 13
Set-TextCell $ws 13 4 '1.628.37'
Set-TextCell $ws 13 5 '  -1.27%  '

# Row 14
Set-TextCell $ws 14 4 '0.5483'
Set-TextCell $ws 14 5 '  +0.29%  '

# Row 15
Set-TextCell $ws 15 4 '0.0₅7759'
Set-TextCell $ws 15 5 '  -1.76%  '

# Row 16
Set-TextCell $ws 16 4 '64.36'
Set-TextCell $ws 16 5 '  -0.66%  '

# Row 17
Set-TextCell $ws 17 4 '26.014.35'
Set-TextCell $ws 17 5 '  +0.01%  '

# Row 18
Set-TextCell $ws 18 5 '  -0.40%  '

# Row 19
Set-TextCell $ws 19 4 '196.85'
Set-TextCell $ws 19 5 '  -0.96%  '

# Row 20
Set-TextCell $ws 20 4 '4.436'
Set-TextCell $ws 20 5 '  -0.74%  '

# Row 21
Set-TextCell $ws 21 4 '9.949'
Set-TextCell $ws 21 5 '  -0.76%  '

# Row 22
Set-TextCell $ws 22 4 '6.102'
Set-TextCell $ws 22 5 '  +0.35%  '

# Row 23
Set-TextCell $ws 23 5 '  -0.55%  '

# Row 24
Set-TextCell $ws 24 4 '1.893'
Set-TextCell $ws 24 5 '  +1.71%  '

# Row 25
Set-TextCell $ws 25 4 '143.81'
Set-TextCell $ws 25 5 '  +2.65%  '

# Row 26
Set-TextCell $ws 26 4 '0.1241'
Set-TextCell $ws 26 5 '  +7.85%  '

# Row 27
Set-TextCell $ws 27 4 '6.877'
Set-TextCell $ws 27 5 '  -0.46%  '

# Row 28
Set-TextCell $ws 28 4 '15.64'
Set-TextCell $ws 28 5 '  -0.84%  '

# Row 29
Set-TextCell $ws 29 5 '  -0.31%  '

# Row 30
Set-TextCell $ws 30 4 '0.04888'
Set-TextCell $ws 30 5 '  -2.83%  '

# Row 31
Set-TextCell $ws 31 4 '3.278'
Set-TextCell $ws 31 5 '  -0.45%  '

# Row 32
Set-TextCell $ws 32 4 '3.228'
Set-TextCell $ws 32 5 '  +0.76%  '

# Row 33
Set-TextCell $ws 33 4 '1.546'
Set-TextCell $ws 33 5 '  -0.04%  '

# Row 34
Set-TextCell $ws 34 5 '  +0.53%  '

# Row 35
Set-TextCell $ws 35 4 '0.9167'
Set-TextCell $ws 35 5 '  +2.38%  '

# Row 36
Set-TextCell $ws 36 4 '2.572'
Set-TextCell $ws 36 5 '  -0.69%  '

# Row 37
Set-TextCell $ws 37 4 '0.5560'
Set-TextCell $ws 37 5 '  +0.26%  '

# Row 38
Set-TextCell $ws 38 4 '1.090.71'
Set-TextCell $ws 38 5 '  -4.10%  '

# Row 39
Set-TextCell $ws 39 4 '0.01573'
Set-TextCell $ws 39 5 '  +0.57%  '

# Row 40
Set-TextCell $ws 40 5 '  -0.51%  '

# Row 41
Set-TextCell $ws 41 4 '2.526'
Set-TextCell $ws 41 5 '  -1.31%  '

# Row 42
Set-TextCell $ws 42 4 '5.604'
Set-TextCell $ws 42 5 '  -1.22%  '

# Row 43
Set-TextCell $ws 43 4 '0.8058'
Set-TextCell $ws 43 5 '  -1.48%  '

# Row 44
Set-TextCell $ws 44 4 '99.22'
Set-TextCell $ws 44 5 '  -0.84%  '

# Row 45
Set-TextCell $ws 45 2 'RocketPoolETH'
Set-TextCell $ws 45 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell $ws 45 4 '1.775.99'
Set-TextCell $ws 45 5 '  -0.48%  '

# Row 46
Set-TextCell $ws 46 2 'BabyDogeCoin'
Set-TextCell $ws 46 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws 46 4 '0.0₈118'
Set-TextCell $ws 46 5 '  -6.04%  '

# Row 47
Set-TextCell $ws 47 5 '  +0.07%  '

# Row 48
Set-TextCell $ws 48 4 '55.59'
Set-TextCell $ws 48 5 '  +0.58%  '

# Row 49
Set-TextCell $ws 49 4 '1.005'
Set-TextCell $ws 49 5 '  +0.16%  '

# Row 50
Set-TextCell $ws 50 5 '  +2.52%  '

# Row 51
Set-TextCell $ws 51 4 '7.565'
Set-TextCell $ws 51 5 '  +2.08%  '
